# Update the division problems table per c986bee.
$d = $word.ActiveDocument

# The value "58÷4=" occurs twice in the table, at different cells, with two
# different replacements, so those two are targeted directly by cell
# position (row 9 col 4, and row 17 col 2) rather than by text search.
$t = $d.Tables.Item(1)
$t.Cell(9, 4).Range.Text = "59÷4="
$t.Cell(17, 2).Range.Text = "67÷6="

# All remaining values are unique in the document, so a plain Find/Replace
# is sufficient for each.
$d.Content.Find.Execute("23÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "53÷8=", 2)
$d.Content.Find.Execute("58÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "66÷9=", 2)
$d.Content.Find.Execute("43÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "97÷7=", 2)
$d.Content.Find.Execute("58÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷3=", 2)
$d.Content.Find.Execute("70÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷3=", 2)
$d.Content.Find.Execute("19÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷3=", 2)
$d.Content.Find.Execute("36÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "13÷7=", 2)
$d.Content.Find.Execute("73÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷4=", 2)
$d.Content.Find.Execute("75÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "89÷7=", 2)
$d.Content.Find.Execute("76÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "62÷3=", 2)
$d.Content.Find.Execute("88÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "89÷4=", 2)
$d.Content.Find.Execute("66÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷4=", 2)
$d.Content.Find.Execute("82÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "77÷6=", 2)
$d.Content.Find.Execute("83÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "10÷5=", 2)
$d.Content.Find.Execute("24÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "11÷4=", 2)
$d.Content.Find.Execute("72÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "63÷9=", 2)
$d.Content.Find.Execute("37÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "86÷9=", 2)
$d.Content.Find.Execute("81÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "42÷4=", 2)
$d.Content.Find.Execute("36÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "10÷8=", 2)
$d.Content.Find.Execute("67÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "38÷5=", 2)
$d.Content.Find.Execute("60÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "64÷2=", 2)
$d.Content.Find.Execute("73÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "49÷7=", 2)
$d.Content.Find.Execute("38÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "45÷3=", 2)
